$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New data set (rows 16-23), replacing the previous state:
#  - row 16 now holds CEILY CENETH CESPEDES CORONADO (1143344960), periodo 1912
#  - rows 17-23 now hold LOLYS JOSEFINA MENDOZA MEDINA (22999419), periodos 2101..2107 in ascending order

$ws.Range("C16").Value = "1143344960"
$ws.Range("D16").Value = "CEILY CENETH CESPEDES CORONADO"
$ws.Range("E16").Value = "1912"
$ws.Range("F16").Value = 33125
$ws.Range("G16").Value = 877803

$ws.Range("C17").Value = "22999419"
$ws.Range("D17").Value = "LOLYS JOSEFINA MENDOZA MEDINA"
$ws.Range("E17").Value = "2101"
$ws.Range("F17").Value = 40000
$ws.Range("G17").Value = 1000000

$ws.Range("C18").Value = "22999419"
$ws.Range("D18").Value = "LOLYS JOSEFINA MENDOZA MEDINA"
$ws.Range("E18").Value = "2102"
$ws.Range("F18").Value = 40000
$ws.Range("G18").Value = 1000000

$ws.Range("C19").Value = "22999419"
$ws.Range("D19").Value = "LOLYS JOSEFINA MENDOZA MEDINA"
$ws.Range("E19").Value = "2103"
$ws.Range("F19").Value = 40000
$ws.Range("G19").Value = 1000000

$ws.Range("C20").Value = "22999419"
$ws.Range("D20").Value = "LOLYS JOSEFINA MENDOZA MEDINA"
$ws.Range("E20").Value = "2104"
$ws.Range("F20").Value = 40000
$ws.Range("G20").Value = 1000000

$ws.Range("C21").Value = "22999419"
$ws.Range("D21").Value = "LOLYS JOSEFINA MENDOZA MEDINA"
$ws.Range("E21").Value = "2105"
$ws.Range("F21").Value = 40000
$ws.Range("G21").Value = 1000000

$ws.Range("C22").Value = "22999419"
$ws.Range("D22").Value = "LOLYS JOSEFINA MENDOZA MEDINA"
$ws.Range("E22").Value = "2106"
$ws.Range("F22").Value = 40000
$ws.Range("G22").Value = 1000000

$ws.Range("C23").Value = "22999419"
$ws.Range("D23").Value = "LOLYS JOSEFINA MENDOZA MEDINA"
$ws.Range("E23").Value = "2107"
$ws.Range("F23").Value = 33333
$ws.Range("G23").Value = 1000000

$wb.Save()
